$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# text (NumberFormat "@") before assignment, otherwise Excel coerces them
# to floating point values instead of preserving the literal text.
$textForceCells = @("D5", "D7", "D9", "D10", "D12", "D14", "D19", "D20", "D22", "D26", "D29", "D36", "D37", "D38", "D41", "D44", "D45", "D47", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price / volume(1h) values scraped by the Action run.
$ws.Range("D2").Value = "38.024.43"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "2.048.66"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "227.85"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").Value = "60.98"
$ws.Range("E7").Value = "  +8.19%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.383"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "0.0817"
$ws.Range("E10").Value = "  +2.08%  "
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "14.72"
$ws.Range("E12").Value = "  +1.73%  "
$ws.Range("D13").Value = "2.352.40"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").Value = "21.08"
$ws.Range("E14").Value = "  +3.40%  "
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("D17").Value = "2.044.38"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "37.985.75"
$ws.Range("E18").Value = "  +1.56%  "
$ws.Range("D19").Value = "6.14"
$ws.Range("E19").Value = "  -1.62%  "
$ws.Range("D20").Value = "69.70"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").Value = "0.0₃0827"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "224.55"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("D26").Value = "166.55"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("D29").Value = "18.93"
$ws.Range("E29").Value = "  +0.82%  "
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("D36").Value = "6.30"
$ws.Range("E36").Value = "  +9.28%  "
$ws.Range("D37").Value = "2.29"
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("D38").Value = "3.24"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "1.537.61"
$ws.Range("E40").Value = "  +4.21%  "
$ws.Range("D41").Value = "97.57"
$ws.Range("E41").Value = "  +2.72%  "
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("D44").Value = "16.61"
$ws.Range("E44").Value = "  +1.61%  "
$ws.Range("D45").Value = "0.0927"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "4.01"
$ws.Range("E47").Value = "  -5.60%  "
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").Value = "7.07"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").Value = "2.241.55"
$ws.Range("E51").Value = "  +0.77%  "

# Restore the default cell style so the forced text-format cells match
# the original (unstyled) look of the rest of the table.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
